$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row -> (DAMSLTag, DialogAct)
$updates = @(
    @{ Row = 12; Tag = "%";   Act = "Uninterpretable" },
    @{ Row = 28; Tag = "qy";  Act = "Yes-No-Question" },
    @{ Row = 32; Tag = "sd";  Act = "Statement-non-opinion" },
    @{ Row = 54; Tag = "aa";  Act = "Agree/Accept" },
    @{ Row = 60; Tag = "sv";  Act = "Statement-opinion" },
    @{ Row = 69; Tag = "aa";  Act = "Agree/Accept" },
    @{ Row = 76; Tag = "sd";  Act = "Statement-non-opinion" },
    @{ Row = 81; Tag = "sv";  Act = "Statement-opinion" },
    @{ Row = 87; Tag = "ba";  Act = "Appreciation" },
    @{ Row = 96; Tag = "aa";  Act = "Agree/Accept" }
)

foreach ($u in $updates) {
    $ws.Cells.Item($u.Row, 9).Value = $u.Tag
    $ws.Cells.Item($u.Row, 10).Value = $u.Act
}
